$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Debbie Gray) - clear first/last name, add office/phone/email
$ws.Range("C2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("G2").Value = "Ag Science, Room 310"
$ws.Range("H2").Value = "Contact via teams"
$ws.Range("I2").Value = "dgray@uidaho.edu"

# Row 3: Soren Newman -> Harpreet Kaur, Statistician, email
$ws.Range("A3").Value = "https://www.uidaho.edu/cals/people/harpreet-kaur"
$ws.Range("C3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "Statistician"
$ws.Range("I3").Value = "hkaur@uidaho.edu"

# Row 4: Harpreet Kaur -> Madison McGuire, Administrative Specialist, office/phone/email
$ws.Range("A4").Value = "https://www.uidaho.edu/cals/people/madison-mcguire"
$ws.Range("C4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "Administrative Specialist"
$ws.Range("G4").Value = "Ag Science, Room 52 (Dean's Suite)"
$ws.Range("H4").Value = "208-885-6681"
$ws.Range("I4").Value = "madisonmcguire@uidaho.edu"

# Row 5: Madison McGuire -> Savanah Nunes, Media and Communications Manager, phone/email
$ws.Range("A5").Value = "https://www.uidaho.edu/cals/people/savanah-nunes"
$ws.Range("C5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "Media and Communications Manager"
$ws.Range("H5").Value = "208-539-7490"
$ws.Range("I5").Value = "snunes@uidaho.edu"

# Row 6: Savanah Nunes -> Brandi Chastain, clear title
$ws.Range("A6").Value = "https://www.uidaho.edu/cals/people/brandi-chastain"
$ws.Range("C6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# Row 7: Steve Greene -> Brian Kelly, clear title
$ws.Range("A7").Value = "https://www.uidaho.edu/cals/people/brian-kelly"
$ws.Range("C7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

# Row 8: Amy Calabretta -> Carly Schoepflin, Director of Comm & Strategic Initiatives, office/phone/email
$ws.Range("A8").Value = "https://www.uidaho.edu/cals/people/carly-schoepflin"
$ws.Range("C8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "Director of Communications & Strategic Initiatives"
$ws.Range("G8").Value = "Ag Science, Room 58"
$ws.Range("H8").Value = "208-885-4037"
$ws.Range("I8").Value = "craska@uidaho.edu"

# Remove rows 9-11 (Angela Harley, Ann Barrington, Ashley Baker)
$ws.Rows("9:11").Delete()
